$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete obsolete rows from the bottom up so earlier row numbers stay stable
# while each delete is applied.
$ws.Range("60:64").Delete()   # stray trailing rows (blank + leftover record rows)
$ws.Range("10:59").Delete()   # collapse the big empty gap between row 9 and row 60
$ws.Rows(8).Delete()          # duplicate "jisha / daniel / it / tester" row
$ws.Rows(5).Delete()          # "jisha / chauhan / it / tester" row

# Header: "ID No" -> "ID"
$ws.Range("A1").Value = "ID"

# New trailing record that replaces the old stray rows
$ws.Range("A8").Value = " "
$ws.Range("B8").Value = "raj"
$ws.Range("C8").Value = "nair"
$ws.Range("D8").Value = "it"
$ws.Range("E8").Value = "tester"
